$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61
$ws.Range("A61").Value = '2022-03-22 13:07:04'
$ws.Range("B61").Value = "'"
$ws.Range("C61").Value = 0.0
$ws.Range("D61").Value = 'NONE'
$ws.Range("E61").Value = 'NONE'
$ws.Range("F61").Value = 'CMNET'
$ws.Range("G61").Value = 223.0
$ws.Range("H61").Value = 5.0
$ws.Range("I61").Value = 5.0
$ws.Range("J61").Value = 5.0
$ws.Range("K61").Value = 0.0
$ws.Range("L61").Value = 0.0
$ws.Range("M61").Value = 0.0
$ws.Range("N61").Value = 0.0
$ws.Range("O61").Value = 0.0
$ws.Range("P61").Value = 1.0
$ws.Range("Q61").Value = '10F872226797'
$ws.Range("R61").Value = 0.0
$ws.Range("S61").Value = 0.0
$ws.Range("T61").Value = 0.0
$ws.Range("U61").Value = 0.0
$ws.Range("V61").Value = 0.0
$ws.Range("W61").Value = "'"
$ws.Range("X61").Value = 0.0
$ws.Range("Y61").Value = 0.0
$ws.Range("Z61").Value = 0.0
$ws.Range("AA61").Value = 72.0
$ws.Range("AB61").Value = 77.0
$ws.Range("AC61").Value = 114.0
$ws.Range("AD61").Value = 34.0
$ws.Range("AE61").Value = 103.0
$ws.Range("AF61").Value = 151.0
$ws.Range("AG61").Value = '1.0.1'
$ws.Range("AH61").Value = '1.0.1'
$ws.Range("AI61").Value = 1.0
$ws.Range("AJ61").Value = 3600.0
$ws.Range("AK61").Value = 'dataeu.hoymiles.com'
$ws.Range("AL61").Value = 1.0
$ws.Range("AM61").Value = 10081.0
$ws.Range("AN61").Value = 0.0
$ws.Range("AO61").Value = 0.0
$ws.Range("AP61").Value = 0.0
$ws.Range("AQ61").Value = 0.0
$ws.Range("AR61").Value = 0.0
$ws.Range("AS61").Value = 0.0
$ws.Range("AT61").Value = 1647950826
$ws.Range("AU61").Value = 0.0
$ws.Range("AV61").Value = 0.0
$ws.Range("AW61").Value = 0.0
$ws.Range("AX61").Value = 0.0
$ws.Range("AY61").Value = 0.0
$ws.Range("AZ61").Value = 0.0
$ws.Range("BA61").Value = 0.0
$ws.Range("BB61").Value = 0.0
$ws.Range("BC61").Value = 0.0
$ws.Range("BD61").Value = 0.0
$ws.Range("BE61").Value = '0negawsklov0negawsklov'
$ws.Range("BF61").Value = 58
$ws.Range("BG61").Value = 'HomeSweetHome'
$ws.Range("BH61").Value = 0.0
$ws.Range("BI61").Value = 0.0

# Row 62
$ws.Range("A62").Value = '2022-03-22 20:19:27'
$ws.Range("B62").Value = "'"
$ws.Range("C62").Value = 0.0
$ws.Range("D62").Value = 'NONE'
$ws.Range("E62").Value = 'NONE'
$ws.Range("F62").Value = 'CMNET'
$ws.Range("G62").Value = 223.0
$ws.Range("H62").Value = 5.0
$ws.Range("I62").Value = 5.0
$ws.Range("J62").Value = 5.0
$ws.Range("K62").Value = 0.0
$ws.Range("L62").Value = 0.0
$ws.Range("M62").Value = 0.0
$ws.Range("N62").Value = 0.0
$ws.Range("O62").Value = 0.0
$ws.Range("P62").Value = 1.0
$ws.Range("Q62").Value = '10F872226797'
$ws.Range("R62").Value = 0.0
$ws.Range("S62").Value = 0.0
$ws.Range("T62").Value = 0.0
$ws.Range("U62").Value = 0.0
$ws.Range("V62").Value = 0.0
$ws.Range("W62").Value = "'"
$ws.Range("X62").Value = 0.0
$ws.Range("Y62").Value = 0.0
$ws.Range("Z62").Value = 0.0
$ws.Range("AA62").Value = 72.0
$ws.Range("AB62").Value = 77.0
$ws.Range("AC62").Value = 114.0
$ws.Range("AD62").Value = 34.0
$ws.Range("AE62").Value = 103.0
$ws.Range("AF62").Value = 151.0
$ws.Range("AG62").Value = '1.0.1'
$ws.Range("AH62").Value = '1.0.1'
$ws.Range("AI62").Value = 1.0
$ws.Range("AJ62").Value = 3600.0
$ws.Range("AK62").Value = 'dataeu.hoymiles.com'
$ws.Range("AL62").Value = 1.0
$ws.Range("AM62").Value = 10081.0
$ws.Range("AN62").Value = 0.0
$ws.Range("AO62").Value = 0.0
$ws.Range("AP62").Value = 0.0
$ws.Range("AQ62").Value = 0.0
$ws.Range("AR62").Value = 0.0
$ws.Range("AS62").Value = 0.0
$ws.Range("AT62").Value = 1647976771
$ws.Range("AU62").Value = 0.0
$ws.Range("AV62").Value = 0.0
$ws.Range("AW62").Value = 0.0
$ws.Range("AX62").Value = 0.0
$ws.Range("AY62").Value = 0.0
$ws.Range("AZ62").Value = 0.0
$ws.Range("BA62").Value = 0.0
$ws.Range("BB62").Value = 0.0
$ws.Range("BC62").Value = 0.0
$ws.Range("BD62").Value = 0.0
$ws.Range("BE62").Value = '0negawsklov0negawsklov'
$ws.Range("BF62").Value = 58
$ws.Range("BG62").Value = 'HomeSweetHome'
$ws.Range("BH62").Value = 0.0
$ws.Range("BI62").Value = 0.0

# Row 63
$ws.Range("A63").Value = '2022-03-23 23:36:21'
$ws.Range("B63").Value = "'"
$ws.Range("C63").Value = 0.0
$ws.Range("D63").Value = 'NONE'
$ws.Range("E63").Value = 'NONE'
$ws.Range("F63").Value = 'CMNET'
$ws.Range("G63").Value = 223.0
$ws.Range("H63").Value = 5.0
$ws.Range("I63").Value = 5.0
$ws.Range("J63").Value = 5.0
$ws.Range("K63").Value = 0.0
$ws.Range("L63").Value = 0.0
$ws.Range("M63").Value = 0.0
$ws.Range("N63").Value = 0.0
$ws.Range("O63").Value = 0.0
$ws.Range("P63").Value = 1.0
$ws.Range("Q63").Value = '10F872226797'
$ws.Range("R63").Value = 0.0
$ws.Range("S63").Value = 0.0
$ws.Range("T63").Value = 0.0
$ws.Range("U63").Value = 0.0
$ws.Range("V63").Value = 0.0
$ws.Range("W63").Value = "'"
$ws.Range("X63").Value = 0.0
$ws.Range("Y63").Value = 0.0
$ws.Range("Z63").Value = 0.0
$ws.Range("AA63").Value = 72.0
$ws.Range("AB63").Value = 77.0
$ws.Range("AC63").Value = 114.0
$ws.Range("AD63").Value = 34.0
$ws.Range("AE63").Value = 103.0
$ws.Range("AF63").Value = 151.0
$ws.Range("AG63").Value = '1.0.1'
$ws.Range("AH63").Value = '1.0.1'
$ws.Range("AI63").Value = 1.0
$ws.Range("AJ63").Value = 3600.0
$ws.Range("AK63").Value = 'dataeu.hoymiles.com'
$ws.Range("AL63").Value = 1.0
$ws.Range("AM63").Value = 10081.0
$ws.Range("AN63").Value = 0.0
$ws.Range("AO63").Value = 0.0
$ws.Range("AP63").Value = 0.0
$ws.Range("AQ63").Value = 0.0
$ws.Range("AR63").Value = 0.0
$ws.Range("AS63").Value = 0.0
$ws.Range("AT63").Value = 1648074986
$ws.Range("AU63").Value = 0.0
$ws.Range("AV63").Value = 0.0
$ws.Range("AW63").Value = 0.0
$ws.Range("AX63").Value = 0.0
$ws.Range("AY63").Value = 0.0
$ws.Range("AZ63").Value = 0.0
$ws.Range("BA63").Value = 0.0
$ws.Range("BB63").Value = 0.0
$ws.Range("BC63").Value = 0.0
$ws.Range("BD63").Value = 0.0
$ws.Range("BE63").Value = '0negawsklov0negawsklov'
$ws.Range("BF63").Value = 52
$ws.Range("BG63").Value = 'HomeSweetHome'
$ws.Range("BH63").Value = 0.0
$ws.Range("BI63").Value = 0.0

# Strip the quote-prefix marker style picked up from the apostrophe trick
# above so these cells end up with the default style, like the rest of the sheet.
$ws.Range("B61").Style = "Normal"
$ws.Range("W61").Style = "Normal"
$ws.Range("B62").Style = "Normal"
$ws.Range("W62").Style = "Normal"
$ws.Range("B63").Style = "Normal"
$ws.Range("W63").Style = "Normal"
